# Apply Thu Jul 11 11:24:03 UTC 2024 cryptos-list update (GitHub Actions bot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.491.72"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "'3.142.86"
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'534.44"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").Value = "'143.37"
$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'3.143.31"
$ws.Range("E8").Value = "  +1.44%  "

$ws.Range("D9").Value = "'0.450"
$ws.Range("E9").Value = "  +1.90%  "

$ws.Range("D10").Value = "'7.18"
$ws.Range("E10").Value = "  -2.15%  "

$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("D12").Value = "'0.395"
$ws.Range("E12").Value = "  +3.09%  "

$ws.Range("D13").Value = "'3.684.12"
$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("E14").Value = "  +3.35%  "

$ws.Range("D15").Value = "'25.70"
$ws.Range("E15").Value = "  -4.07%  "

$ws.Range("D16").Value = "'0.0000167"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "'58.540.44"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").Value = "'3.141.85"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("D19").Value = "'6.12"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "'12.91"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").Value = "'343.13"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").Value = "'0.514"
$ws.Range("E24").Value = "  +1.90%  "

$ws.Range("D25").Value = "'67.82"

$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "'0.0₃0934"
$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("D29").Value = "'7.54"
$ws.Range("E29").Value = "  +4.19%  "

$ws.Range("D30").Value = "'6.48"
$ws.Range("E30").Value = "  -2.79%  "

$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("D32").Value = "'1.90"
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").Value = "'4.81"
$ws.Range("E35").Value = "  +3.76%  "

$ws.Range("D36").Value = "'158.17"
$ws.Range("E36").Value = "  +2.50%  "

$ws.Range("E37").Value = "  +3.63%  "

$ws.Range("D38").Value = "'26.29"
$ws.Range("E38").Value = "  -2.09%  "

$ws.Range("E39").Value = "  -3.12%  "

$ws.Range("E40").Value = "  +11.41%  "

$ws.Range("D41").Value = "'0.0673"
$ws.Range("E41").Value = "  -0.70%  "

$ws.Range("D42").Value = "'0.711"
$ws.Range("E42").Value = "  +4.84%  "

$ws.Range("D43").Value = "'4.02"
$ws.Range("E43").Value = "  +3.74%  "

$ws.Range("D44").Value = "'3.184.50"
$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("D45").Value = "'36.62"
$ws.Range("E45").Value = "  -0.52%  "

$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "'0.0266"
$ws.Range("E47").Value = "  +3.28%  "

$ws.Range("D48").Value = "'2.303.15"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("D49").Value = "'1.01"
$ws.Range("E49").Value = "  +5.04%  "

$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("E51").Value = "  +2.05%  "

